# Add new row 11 ("2021年") to Sheet1, matching the 2012-2020 rows already present.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the year label cell (A10) onto the new label cell (A11)
# so the new row keeps the same bold/centered/bordered formatting as the others.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 66.39
$ws.Range("C11").Value = 7.06
$ws.Range("D11").Value = 0.5600000000000001
$ws.Range("E11").Value = ""  # no data for this indicator in 2021
$ws.Range("F11").Value = 40.22
$ws.Range("G11").Value = 304.87
$ws.Range("H11").Value = 27.52
$ws.Range("I11").Value = 14.37
$ws.Range("J11").Value = -0.59
$ws.Range("K11").Value = 4891.75
$ws.Range("L11").Value = 0.26
$ws.Range("M11").Value = 9.06
$ws.Range("N11").Value = 1.5
$ws.Range("O11").Value = 3.04
$ws.Range("P11").Value = 231.9
$ws.Range("Q11").Value = 34.77
$ws.Range("R11").Value = 3.98
$ws.Range("S11").Value = 10.93
$ws.Range("T11").Value = 148.83
$ws.Range("U11").Value = -120.49
$ws.Range("V11").Value = -85.56999999999999
$ws.Range("W11").Value = 730.87
$ws.Range("X11").Value = 59.47
$ws.Range("Y11").Value = 2411.93
$ws.Range("Z11").Value = 33.85
$ws.Range("AA11").Value = 0.13
$ws.Range("AB11").Value = 139.45
$ws.Range("AC11").Value = 141.67
$ws.Range("AD11").Value = 13.35
$ws.Range("AE11").Value = -0.75
$ws.Range("AF11").Value = 98.17
$ws.Range("AG11").Value = 9.199999999999999
$ws.Range("AH11").Value = 16.69
$ws.Range("AI11").Value = -25.74
$ws.Range("AJ11").Value = 5.61
$ws.Range("AK11").Value = 27.46
$ws.Range("AL11").Value = 15.72
$ws.Range("AM11").Value = 80.47
$ws.Range("AN11").Value = 29.31
$ws.Range("AO11").Value = 12.1
$ws.Range("AP11").Value = 272.29
$ws.Range("AQ11").Value = 121.8
